# Update column F (dSF) values for specific rows to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -4
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -5
$ws.Range("F18").Value = -8
$ws.Range("F19").Value = -6
$ws.Range("F24").Value = -1
$ws.Range("F25").Value = -2
